$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-40, 44-51: update Price (D) and Volume(1h) (E) values.
# The Price column (D) is forced back to Text format before each
# write so numeric-looking strings (e.g. "242.36") are not silently
# auto-converted to numbers by Excel, matching the source data which
# stores these as inline strings.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.243.21'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.860.92'
$ws.Range("E3").Value = '  -1.09%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.36'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9996'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3131'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07803'
$ws.Range("E9").Value = '  -2.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.26'
$ws.Range("E10").Value = '  -4.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07999'
$ws.Range("E11").Value = '  -4.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.864.50'
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '94.08'
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.181'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6973'
$ws.Range("E15").Value = '  -3.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.396'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.251.06'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008294'
$ws.Range("E18").Value = '  -3.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '253.22'
$ws.Range("E19").Value = '  +4.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.13'
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.108.30'
$ws.Range("E21").Value = '  -1.83%  '
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.535'
$ws.Range("E23").Value = '  -4.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9995'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1565'
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.996'
$ws.Range("E26").Value = '  -0.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.82'
$ws.Range("E27").Value = '  -2.31%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.90'
$ws.Range("E28").Value = '  +1.46%  '
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.310'
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.267'
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.209'
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05275'
$ws.Range("E33").Value = '  -2.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.888'
$ws.Range("E34").Value = '  -3.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7486'
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.156'
$ws.Range("E36").Value = '  -2.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.709'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01870'
$ws.Range("E38").Value = '  -1.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.245.07'
$ws.Range("E39").Value = '  -3.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.734'
$ws.Range("E40").Value = '  -0.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.54'
$ws.Range("E44").Value = '  -5.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9987'
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.005.87'
$ws.Range("E47").Value = '  -1.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5189'
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.785'
$ws.Range("E49").Value = '  -1.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.479'
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4299'
$ws.Range("E51").Value = '  -2.19%  '

# Rows 41-43: reorder rows (rotate B/C/D/E values)
$ws.Range("B41").Value = 'Quant'
$ws.Range("C41").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '111.01'
$ws.Range("E41").Value = '  -0.71%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.129'
$ws.Range("E42").Value = '  -7.14%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8951'
$ws.Range("E43").Value = '  -2.38%  '
